# Auto bat commit friday
#
# 1. On the existing "booking_data" sheet, the user selected the whole
#    sheet (Ctrl+A) - active cell ends up "C9" with the selection
#    spanning the entire grid (A1:XFD1048576), and it's no longer the
#    sheet shown when the file is reopened.
# 2. A brand-new sheet called "Book" is appended after "booking_data"
#    and filled with a small 2-row/6-column "books" dataset. This sheet
#    becomes the active / selected tab, with cell E11 selected.

$wb = $excel.ActiveWorkbook

$bookingData = $wb.Worksheets.Item("booking_data")
$bookingData.Activate()
$bookingData.Cells.Select()

# Add the new "Book" sheet as the very last tab in the workbook.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Book"

# Borrow the existing look (bold/shaded header row, bordered data rows,
# text-formatted "pageCount" column, date-style "excerpt"/"publishDate"
# columns) from the "booking_data" sheet so the new sheet's styling
# matches the rest of the workbook instead of inventing new styles.
$bookingData.Range("A1:F1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)
$bookingData.Range("A2:F2").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$ws.Range("A3:F3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row.
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "title"
$ws.Range("C1").Value = "description"
$ws.Range("D1").Value = "pageCount"
$ws.Range("E1").Value = "excerpt"
$ws.Range("F1").Value = "publishDate"

# Row 2.
$ws.Range("F2").Value = "2023-01-02T16:51:12.898Z"
$ws.Range("B2").Value = "Sherlok"
$ws.Range("C2").Value = "The Sign of Four"
$ws.Range("D2").Value = "300"
$ws.Range("E2").Value = "Thriller"
$ws.Range("A2").Value = 983242

# Row 3.
$ws.Range("B3").Value = "Romana"
$ws.Range("C3").Value = "The Dangerous Path"
$ws.Range("D3").Value = "400"
$ws.Range("E3").Value = "SiFi"
$ws.Range("F3").Value = "2023-03-03T16:51:12.898Z"
$ws.Range("A3").Value = 879374

$ws.Columns("B:F").AutoFit()

$ws.Range("E11").Select()
